$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F (shifts old F "District" -> G)
$ws.Columns.Item(6).Insert()

# Header
$ws.Range("F2").Value = "Address"

# Address values per row (new column F)
$ws.Range("F3").Value = "Govt. High School HudukulaBangarpet"
$ws.Range("F4").Value = "G H S Harati"
$ws.Range("F5").Value = "A H S RayalamanadinneMulbagal"
$ws.Range("F6").Value = "G J C Gownipalli (V)Sreenivasapura"
$ws.Range("F7").Value = "B H S High School AbbenahalliMalur"
$ws.Range("F8").Value = "Govt. Junior CollegeTayalurMulbagal"
$ws.Range("F9").Value = "Govt. High School KyasamballiBangarpet"
$ws.Range("F10").Value = "Govt Jr College M N HalliMulbagal"
$ws.Range("F11").Value = "S V V High School NangaliMulbagal"
$ws.Range("F12").Value = "Citizen High School ThayalurMulbagal"
$ws.Range("F13").Value = "Field Marshal K M Cariyappa High School HanchalagateBangarpet"
$ws.Range("F14").Value = "Govt. High School Alangur Cross Mulbagal"
$ws.Range("F15").Value = ""
$ws.Range("F16").Value = "G J C Tekal Malur"
$ws.Range("F17").Value = "Vivekananda Rural High SchoolN G HulkurBangarpet"
$ws.Range("F18").Value = "Govt. Junior College Vemgal"
$ws.Range("F19").Value = "G H S DevarayasamudraMulbagal"
$ws.Range("F20").Value = "G H S RonurSrinivasapur"
$ws.Range("F21").Value = "G J C DoddanayakanahalliMalur"
$ws.Range("F22").Value = "G J C D N DoddiMalur"
$ws.Range("F23").Value = "G J C (H S) MasthiMalur"
$ws.Range("F24").Value = ""
$ws.Range("F25").Value = "Mysore mine High SchoolK G F BlockBangarpet"
$ws.Range("F26").Value = "Nethaji High School MustoorMulbagal"
$ws.Range("F27").Value = "G H S RoyalpadSreenivasapura"
$ws.Range("F28").Value = "Sabaramathi High SchoolSugatur"
$ws.Range("F29").Value = "Adarsha Vidyalaya (RMSA) Bangarpet"
$ws.Range("F30").Value = "G H S KesaragereMalur"
$ws.Range("F31").Value = "Govt. High School S MadamangalaBangarpet"
$ws.Range("F32").Value = "G J C Boys Malur"
$ws.Range("F33").Value = "Anjandri High School Emmenatha Mulbagal"
$ws.Range("F34").Value = "Robertsonpet K G "
$ws.Range("F35").Value = "G H S ChikkathirupathiMalur"
$ws.Range("F36").Value = "G H S NG HulkurBangarpet"
$ws.Range("F37").Value = "J H S Kembodi"
$ws.Range("F38").Value = "Govt. Jr. College for BoysSrinivasapura"
$ws.Range("F39").Value = "M E M High School"
$ws.Range("F40").Value = "G B J C Bangarapet"
$ws.Range("F41").Value = "G H S Kyalanure"
$ws.Range("F42").Value = "G H P S Shettikothanur"
$ws.Range("F43").Value = "J V T H S PathimitteMulbagal"
$ws.Range("F44").Value = "G H S BudikoteBangarpet"
$ws.Range("F45").Value = "Govt. High School Thippadoddi Mulbagal"
$ws.Range("F46").Value = "Mathru Bhoomi High School Kodikannur"
$ws.Range("F47").Value = "S K R S High School MitturMulbagal"
$ws.Range("F48").Value = "G H S VirupakshiMulbagal"
$ws.Range("F49").Value = "G H S Annihalli"
$ws.Range("F50").Value = "G J C KamasamudraK G F RangeBangarpet"
$ws.Range("F51").Value = "G H S Gonamakanahalli"
$ws.Range("F52").Value = "Govt. High School Kamadhenuhalli"
$ws.Range("F53").Value = "G H SGullahalliBangarpet"

$ws.Range("A1").Select()
